$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Img Url" column (H) ------------------------------------------------
$ws.Cells.Item(1, 8).Value = "Img Url"

$ws.Cells.Item(23, 8).Value = "https://www.makro.co.za/sys-master/images/h06/h98/9800090222622/silo-MIN_237823_EAA_large"
$ws.Cells.Item(24, 8).Value = "https://www.builders.co.za/_ui/responsive/theme-yellow/images/products/product-image-transparent.png"
$ws.Cells.Item(25, 8).Value = "https://www.makro.co.za/sys-master/images/hc0/h4f/9800097300510/silo-MIN_159320_EAA_large"
$ws.Cells.Item(26, 8).Value = "https://www.makro.co.za/sys-master/images/hba/h86/9800090746910/silo-MIN_285367_EAA_large"
$ws.Cells.Item(27, 8).Value = "https://www.makro.co.za/sys-master/images/hba/h86/9800090746910/silo-MIN_285367_EAA_large"
$ws.Cells.Item(28, 8).Value = "https://3pmedia.leroymerlin.co.za/SOURCE/a424840859ea40f38d972c635c8539ec"
$ws.Cells.Item(29, 8).Value = "https://3pmedia.leroymerlin.co.za/SOURCE/10f337c37a7845719d921542d1415339"
$ws.Cells.Item(30, 8).Value = "https://www.makro.co.za/sys-master/images/hba/h86/9800090746910/silo-MIN_285367_EAA_large"
$ws.Cells.Item(31, 8).Value = "https://www.makro.co.za/sys-master/images/hba/h86/9800090746910/silo-MIN_285367_EAA_large"

# Widen the columns so the (now longer) urls are readable.
$ws.Columns(5).ColumnWidth = 44.21875
$ws.Columns(8).ColumnWidth = 91.88671875

# --- Filter the product list down to Generators only ----------------------------
$ws.Range("A1:G72").AutoFilter(1, @("Generator"), 7)

# Record the filter database name (mirrors what Excel stores when AutoFilter is applied).
$fname = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$72")
$fname.Visible = $false

# --- Update cursor / selection ---------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("H73").Select()
